# edit.ps1 - Apply GearCardData.xlsx changes described in the commit:
# "Now with Mimic Tooth the player can place at most 1 card in the
#  opponent side for each row."
#
# Concretely (for this workbook):
#   1. Row 8 (宝箱怪的牙 / Mimic Tooth), column C (effect) text is updated
#      from "向战场出牌时，可以将牌出在敌对侧。"
#      to   "向战场出牌时，每1行可以将至多1张牌出在敌对侧。"
#   2. The cell D11 (imageFile of 战旗/Warbanner) loses its extra
#      "apply fill" style variant and is normalized to the same style
#      used by the rest of column D ("常规 2" without the extra
#      applyFill flag).
#   3. The active selection on the sheet moves from D12 to C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the effect text for the Mimic Tooth card (row 8, column C).
$ws.Range("C8").Value = "向战场出牌时，每1行可以将至多1张牌出在敌对侧。"

# 2. Normalize D11's style to match the rest of column D (drop the
#    redundant applyFill-only style variant).
$ws.Range("D11").Style = "常规 2"

# 3. Update the sheet's current selection to C9.
[void]$ws.Range("C9").Select()
